# Auto-update predictions and index for 2025-10-30
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 1.FC Köln - Bayern Munich (final score 1:4, Bayern Munich pick correct)
$ws.Range("A2").Value = "1.FC Köln - Bayern Munich ✓: 1:4"
$ws.Range("G2").Value = "✓"

# Row 3: Arsenal FC - Brighton & Hove Albion (final score 2:0, Arsenal FC pick correct)
$ws.Range("A3").Value = "Arsenal FC ✓ - Brighton & Hove Albion: 2:0"
$ws.Range("G3").Value = "✓"

# Row 4: Inter Milan - ACF Fiorentina (final score 3:0, Inter Milan pick correct)
$ws.Range("A4").Value = "Inter Milan ✓ - ACF Fiorentina: 3:0"
$ws.Range("C4").Value = 71
$ws.Range("G4").Value = "✓"

# Row 6: Celtic FC - Falkirk FC (final score 4:0, Celtic FC pick correct)
$ws.Range("A6").Value = "Celtic FC ✓ - Falkirk FC: 4:0"
$ws.Range("C6").Value = 66
$ws.Range("G6").Value = "✓"

# Row 7: new fixture appended - AS Roma - Parma Calcio 1913 (final score 2:1, AS Roma pick correct)
$ws.Range("A7").Value = "AS Roma ✓ - Parma Calcio 1913: 2:1"
$ws.Range("B7").Value = "AS Roma"
$ws.Range("C7").Value = 55
$ws.Range("D7").Value = 86
$ws.Range("E7").Value = 91
$ws.Range("F7").Value = 1.48
$ws.Range("G7").Value = "✓"
